$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the filter/search labels in column E (shared-string text changes)
$ws.Range("E2").Value = "?"
$ws.Range("E10").Value = "Pa monitorim"

# Update the visible view state: scroll so row 3 is the top row and
# select E11 (matches the sheetView/selection change in the diff)
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 3
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("E11").Select()
